# Docx writer: Use different style for block quotes in notes.
#
# Adds a new paragraph style "Footnote Block Text" (styleId
# "FootnoteBlockText"), based on "Footnote Text" (next style also
# "Footnote Text"), with the same uiPriority/unhideWhenUsed/qFormat
# flags and the same block-quote spacing/indent as the existing
# "Block Text" style:
#   spacing: after=100 twips (5pt), before=100 twips (5pt)
#   indent : firstLine=0, left=480 twips (24pt), right=480 twips (24pt)

$d = $word.ActiveDocument

$footnoteText = $d.Styles("Footnote Text")

$blockStyle = $d.Styles.Add("Footnote Block Text", 1)

$blockStyle.BaseStyle = $footnoteText
$blockStyle.NextParagraphStyle = $footnoteText
$blockStyle.Priority = 9
$blockStyle.UnhideWhenUsed = $true
$blockStyle.QuickStyle = $true

$blockStyle.ParagraphFormat.SpaceAfter = 5
$blockStyle.ParagraphFormat.SpaceBefore = 5
$blockStyle.ParagraphFormat.FirstLineIndent = 0
$blockStyle.ParagraphFormat.LeftIndent = 24
$blockStyle.ParagraphFormat.RightIndent = 24
